$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preview import teacher: clear the stale "maso" (B7) and "hodem" (C7) values
# that were left over in row 7, keeping only the "ten" (D7) value "Thanh Bình".
$ws.Range("B7:C7").Clear()

# Move the active selection to C7, matching the saved cursor position.
$ws.Range("C7").Select()
